# Retraining the forecast model for Dragosel Laslea:
# Shift all timestamps in column A (rows 2-97) forward by 2 days,
# and update the forecasted production values in column B for the
# ramp-up period (rows 29-42) to reflect the retrained model output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every timestamp in column A (rows 2 through 97) forward by 2 days,
# keeping the existing time-of-day fraction intact.
for ($row = 2; $row -le 97; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value2 = $cell.Value2 + 2
}

# Updated forecast values (column B) for rows 29-42, produced by the
# retrained model.
$newValues = @{
    29 = 0
    30 = 14
    32 = 66
    33 = 116
    34 = 162
    35 = 227
    36 = 276
    37 = 361
    38 = 395
    39 = 427
    40 = 487
    41 = 542
    42 = 539
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 2).Value2 = $newValues[$row]
}
